$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zub_Gel")

# Reorder the data rows 3-6: the "Menge / Laenge" row moves up to row 3,
# pushing the "P_Art_S1" row down to row 4; the "Menge / Stueck1" row
# moves to row 5, and the "P_Art_S2" row moves down to row 6.

$ws.Range("A3").Value = "Zahl"
$ws.Range("B3").Value = "Menge / Länge"
$ws.Range("C3").Value = "Menge_L"
$ws.Range("D3").Value = ""

$ws.Range("A4").Value = "Mehrfach"
$ws.Range("B4").Value = "Stk. Artikel"
$ws.Range("C4").Value = "P_Art_S1"
$ws.Range("D4").Value = "Balkonblenden-Halterung für je 1 BV 160-400:26, Balkonblenden-Abwinkelung für BV 160 - 400:60, Eckelement min 90mm, Geländer-Blumenkästen 85cm:148, Geländer-Blumenkästen 115cm:171, Geländer-Blumenkästen 165cm:211, Geländer-Blumenkästen 220cm:246"

$ws.Range("A5").Value = "Zahl"
$ws.Range("B5").Value = "Menge / Stück1"
$ws.Range("C5").Value = "Menge_S1"

$ws.Range("A6").Value = "Mehrfach"
$ws.Range("B6").Value = "Stk. Artikel"
$ws.Range("C6").Value = "P_Art_S2"
$ws.Range("D6").Value = "Balkonblenden-Halterung für je 1 BV 160-400:26, Balkonblenden-Abwinkelung für BV 160 - 400:60, Eckelement min 90mm, Geländer-Blumenkästen 85cm:148, Geländer-Blumenkästen 115cm:171, Geländer-Blumenkästen 165cm:211, Geländer-Blumenkästen 220cm:246"

# Switch the active sheet/tab to "Zub_Gel" and select cell D13, matching
# the saved view state (tabSelected + selection move from Brix_Gel_Stab
# to Zub_Gel).
$ws.Activate()
$ws.Range("D13").Select()
